# Generate Report for Handback
#
# The row for file "03f1368e-e57c-41da-a317-409030d904b9.md" moves from
# "Ready for handoff" to "Handback transform failed" because the handback
# transform detected a filename mismatch. Reflect that on the Overview
# sheet (both locale status columns) and on each locale detail sheet's
# Status column, and record the error detail for each locale in the
# "Error Detail" column (K).

$wb = $excel.ActiveWorkbook

$newStatus = "Handback transform failed"

# --- Overview sheet: zh-cn (B3) and de-de (C3) status for the 03f1368e... file
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# --- zh-cn detail sheet
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = $newStatus
$wsZh.Range("K3").Value = "Handback file name: p53oqv03.hin is different with handoff file name: 03f1368e-e57c-41da-a317-409030d904b9.093303959cdf9140d117ca9b01b45c5448ed8cd1.zh-cn."

# --- de-de detail sheet
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = $newStatus
$wsDe.Range("K3").Value = "Handback file name: p53oqv03.hin is different with handoff file name: 03f1368e-e57c-41da-a317-409030d904b9.093303959cdf9140d117ca9b01b45c5448ed8cd1.de-de."
